$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 14 (pushes old rows 14-18 down to 17-21,
# inheriting row 14's formatting - including the quote-prefix style on column C)
$ws.Range("B14:B16").EntireRow.Insert()

# Populate the new row 14 with the relocated deprecation-check test case
$ws.Range("B14").Value = "arrLength"
$ws.Range("C14").Value = "'= """".bytes.length"

# Update the selected cell to reflect the new layout
$ws.Range("C15").Select()

Write-Output "done"
